# Update cryptocurrency price/volume figures per the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.551.50"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "2.642.80"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'604.15"
$ws.Range("E5").Value = "  +2.15%  "

$ws.Range("D6").Value = "'146.46"
$ws.Range("E6").Value = "  +1.80%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  +0.30%  "

$ws.Range("E9").Value = "  +2.06%  "

$ws.Range("D10").Value = "'5.58"
$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("E11").Value = "  +4.29%  "

$ws.Range("E12").Value = "  -0.26%  "

$ws.Range("D13").Value = "'27.62"
$ws.Range("E13").Value = "  +0.87%  "

$ws.Range("D14").Value = "3.118.72"
$ws.Range("E14").Value = "  -0.12%  "

$ws.Range("D15").Value = "63.406.23"
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("E16").Value = "  +1.73%  "

$ws.Range("D17").Value = "2.631.19"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("D18").Value = "'11.52"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("E19").Value = "  +4.98%  "

$ws.Range("D20").Value = "'344.61"
$ws.Range("E20").Value = "  +1.32%  "

$ws.Range("E21").Value = "  +2.95%  "

$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").Value = "'5.57"
$ws.Range("E23").Value = "  -3.31%  "

$ws.Range("D24").Value = "'66.67"
$ws.Range("E24").Value = "  -0.81%  "

$ws.Range("E25").Value = "  +1.84%  "

$ws.Range("D26").Value = "'586.14"
$ws.Range("E26").Value = "  +8.19%  "

$ws.Range("E27").Value = "  +8.20%  "

$ws.Range("D28").Value = "'1.55"
$ws.Range("E28").Value = "  +1.80%  "

$ws.Range("E29").Value = "  -1.48%  "

$ws.Range("D30").Value = "'8.00"
$ws.Range("E30").Value = "  +2.38%  "

$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("E32").Value = "  +4.02%  "

$ws.Range("E33").Value = "  -2.79%  "

$ws.Range("E34").Value = "  +2.73%  "

$ws.Range("D35").Value = "'5.24"
$ws.Range("E35").Value = "  +7.81%  "

$ws.Range("D36").Value = "'166.95"
$ws.Range("E36").Value = "  -4.52%  "

$ws.Range("E37").Value = "  +1.07%  "

$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("D39").Value = "'1.96"
$ws.Range("E39").Value = "  +8.19%  "

$ws.Range("D40").Value = "'19.12"
$ws.Range("E40").Value = "  +0.27%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").Value = "'168.08"
$ws.Range("E42").Value = "  -2.47%  "

$ws.Range("E43").Value = "  +1.44%  "

$ws.Range("D44").Value = "'22.21"
$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("E45").Value = "  +1.00%  "

$ws.Range("D46").Value = "'0.630"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("E47").Value = "  +3.29%  "

$ws.Range("D48").Value = "'0.0964"
$ws.Range("E48").Value = "  +0.29%  "

$ws.Range("D49").Value = "'1.91"
$ws.Range("E49").Value = "  +12.50%  "

$ws.Range("D50").Value = "'18.83"
$ws.Range("E50").Value = "  +1.02%  "

$ws.Range("E51").Value = "  +1.71%  "
